# Re-derive the workbook from the live Excel application object (the
# pre-bound $wb handed to us loses its .ActiveSheet binding for some
# member accesses in this host, so rebuild it the same way the task
# description does).
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("FLT_EFF_YY")
$ws2 = $wb.Worksheets.Item("FLT_EFF_MM")

# ---------------------------------------------------------------------
# 1. Metadata link text + contact e-mail text (shared on both sheets,
#    cells F1 / F2).
# ---------------------------------------------------------------------
$f1Sheet1 = $ws1.Range("F1")
$f1Sheet1.Value = "Metadata - Single European Sky Portal"

$f2Sheet1 = $ws1.Range("F2")
$f2Sheet1.Value = "pru-support@eurocontrol.int"

$f1Sheet2 = $ws2.Range("F1")
$f1Sheet2.Value = "Metadata - Single European Sky Portal"

$f2Sheet2 = $ws2.Range("F2")
$f2Sheet2.Value = "pru-support@eurocontrol.int"

# ---------------------------------------------------------------------
# 2. Sheet2 no longer carries a mailto: hyperlink on F2 - only the
#    metadata-link hyperlink on F1 remains. Do this before any font
#    tweaks below, because (re)adding a hyperlink re-applies the
#    built-in "Hyperlink" style and would clobber explicit font
#    formatting applied afterwards.
# ---------------------------------------------------------------------
$f2Sheet2.Hyperlinks.Delete()
$ws2.Hyperlinks.Add($f1Sheet2, "http://prudata.webfactional.com/wiki/index.php/Average_horizontal_en-route_inefficiency")

# ---------------------------------------------------------------------
# 3. Font restyle of those same link-styled cells.
# ---------------------------------------------------------------------
$f1Sheet1.Font.Name = "Arial"
$f1Sheet1.Font.Size = 9
$f1Sheet1.Font.Color = 16711680

$f1Sheet2.Font.Name = "Arial"
$f1Sheet2.Font.Size = 9
$f1Sheet2.Font.Color = 13391121

$f2Sheet2.Font.Size = 9

# ---------------------------------------------------------------------
# 4. Outline properties now explicit on both sheets.
# ---------------------------------------------------------------------
$ws1.Outline.SummaryRow = 0
$ws1.Outline.SummaryColumn = 0
$ws2.Outline.SummaryRow = 0
$ws2.Outline.SummaryColumn = 0

# ---------------------------------------------------------------------
# 5. Column widths shrink slightly on both sheets (re-export rescale).
# ---------------------------------------------------------------------
$offset = 0.8333333333333334

$s1c1 = $ws1.Columns.Item(1)
$s1c1.ColumnWidth = 11.5 - $offset
$s1c2 = $ws1.Columns.Item(2)
$s1c2.ColumnWidth = 15.13 - $offset
$s1c3 = $ws1.Columns.Item(3)
$s1c3.ColumnWidth = 9.13 - $offset
$s1c4 = $ws1.Columns.Item(4)
$s1c4.ColumnWidth = 11.75 - $offset
$s1c5 = $ws1.Columns.Item(5)
$s1c5.ColumnWidth = 13.0 - $offset
$s1c6 = $ws1.Columns.Item(6)
$s1c6.ColumnWidth = 11.5 - $offset
$s1c7 = $ws1.Columns.Item(7)
$s1c7.ColumnWidth = 10.38 - $offset

$s2c1 = $ws2.Columns.Item(1)
$s2c1.ColumnWidth = 11.5 - $offset
$s2c2 = $ws2.Columns.Item(2)
$s2c2.ColumnWidth = 15.13 - $offset
$s2c3 = $ws2.Columns.Item(3)
$s2c3.ColumnWidth = 9.13 - $offset
$s2c4 = $ws2.Columns.Item(4)
$s2c4.ColumnWidth = 9.13 - $offset
$s2c5 = $ws2.Columns.Item(5)
$s2c5.ColumnWidth = 17.75 - $offset
$s2c6 = $ws2.Columns.Item(6)
$s2c6.ColumnWidth = 9.0 - $offset

Write-Host "edit applied"
